# Fixed Makefile for automotive/qsort. Added results for automotive/qsort.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New benchmark result rows for automotive/qsort
$ws.Range("A5").Value = "automotive/qsort/qsort_small input_small.dat"
$ws.Range("B5").Value = 0.06
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0.02

$ws.Range("A6").Value = "automotive/qsort/qsort_large input_large.dat"
$ws.Range("B6").Value = 0.29
$ws.Range("C6").Value = 0.05
$ws.Range("D6").Value = 0.12

# Widen column A slightly to fit the new (longer) name strings
$ws.Columns.Item(1).ColumnWidth = 36.8

# Update the active selection like it was left after editing the sheet
[void]$ws.Range("B9").Select()
